# "day 12 vlookup and hlookup.xlsx" — add HLOOKUP formulas to the
# "Sheet2" worksheet (rows 2 and 5, columns C:H), looking up the
# row-1 header labels against the transposed reference table in
# B11:H13. Mirrors what a user gets by typing the formula in C2/C5
# and then dragging the fill handle from D2/D5 across to H2/H5
# (hence the shared-formula group for D:H while C keeps its own
# formula).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Activate()

# Row 2: employee count lookup (row index 2 in the B11:H13 table)
$ws.Range("C2").Formula = "=HLOOKUP(C1,`$B11:`$H13,2,FALSE)"
$ws.Range("D2:H2").Formula = "=HLOOKUP(D1,`$B11:`$H13,2,FALSE)"

# Row 5: address lookup (row index 3 in the B11:H13 table)
$ws.Range("C5").Formula = "=HLOOKUP(C1,`$B11:`$H13,3,FALSE)"
$ws.Range("D5:H5").Formula = "=HLOOKUP(D1,`$B11:`$H13,3,FALSE)"

# Leave the selection where the author ended up after filling the
# formulas across the sheet.
$ws.Range("H9").Select()
